# Add data for 2021-10-17: advance the "through" date from 10-08 to 10-09
# and update the October / Total figures accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab to reflect the new "through" date.
$ws.Name = "Through 2021-10-09"

# Row 11 (September) — only the 2021 totals shift.
$ws.Range("U11").Value = 175
$ws.Range("V11").Value = 0.0223

# Row 12 (October) — label + per-year counts/rates.
$ws.Range("A12").Value = "October (through 10-09)"
$ws.Range("C12").Value = 6
$ws.Range("D12").Value = 0.1429
$ws.Range("F12").Value = 16
$ws.Range("H12").Value = 4
$ws.Range("I12").Value = 13
$ws.Range("J12").Value = 0.2353
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 21
$ws.Range("M12").Value = 0.087
$ws.Range("R12").Value = 38
$ws.Range("U12").Value = 60

# Row 13 (Total) — recomputed grand totals.
$ws.Range("C13").Value = 202
$ws.Range("D13").Value = 0.133
$ws.Range("F13").Value = 399
$ws.Range("G13").Value = 0.1034
$ws.Range("H13").Value = 54
$ws.Range("I13").Value = 590
$ws.Range("J13").Value = 0.0839
$ws.Range("K13").Value = 63
$ws.Range("L13").Value = 508
$ws.Range("M13").Value = 0.1103
$ws.Range("R13").Value = 886
$ws.Range("S13").Value = 0.0564
$ws.Range("U13").Value = 1231
$ws.Range("V13").Value = 0.0603
